$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

$ws.Range("H13").Value = 1566.1666
$ws.Range("I13").Value = 1479.4
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 1479.4
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = -1310.4
$ws.Range("N13").Value = -2338

$ws.Range("H28").Value = 830
$ws.Range("I28").Value = 648
$ws.Range("J28").Value = 1103
$ws.Range("K28").Value = 648
$ws.Range("L28").Value = 1103
$ws.Range("M28").Value = -163

$ws.Range("H100").Value = 2264.8572
$ws.Range("I100").Value = 2264.8572
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2264.8572
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1723.8572

$ws.Range("H113").Value = 8927.083000000001
$ws.Range("I113").Value = 8712
$ws.Range("J113").Value = 10002.5
$ws.Range("K113").Value = 8712
$ws.Range("L113").Value = 10002.5
$ws.Range("M113").Value = -5458

$ws.Range("H138").Value = 4299.6665
$ws.Range("I138").Value = 4439.4
$ws.Range("J138").Value = 4125
$ws.Range("K138").Value = 13318.2
$ws.Range("L138").Value = 12375
$ws.Range("M138").Value = -8178.199999999999
$ws.Range("N138").Value = -22655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3779.6
$ws.Range("I2").Value = 2225
$ws.Range("J2").Value = 9998
$ws.Range("K2").Value = 2225
$ws.Range("L2").Value = 9998
$ws.Range("M2").Value = -2112

$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -385
$ws.Range("N3").Value = -1230

$ws.Range("H8").Value = 10006000
$ws.Range("I8").Value = 20000000
$ws.Range("J8").Value = 12000
$ws.Range("K8").Value = 20000000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = -19999856
$ws.Range("N8").Value = -12288

$ws.Range("H10").Value = 5000500
$ws.Range("I10").Value = 10000000
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 10000000
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -9999830
$ws.Range("N10").Value = -1340

$ws.Range("H11").Value = 1339000
$ws.Range("I11").Value = 2007500
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 2007500
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -2007356
$ws.Range("N11").Value = -2288

$ws.Range("H12").Value = 7700
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -227
$ws.Range("N12").Value = -15346

$ws.Range("H13").Value = 6667331.5
$ws.Range("I13").Value = 20000000
$ws.Range("J13").Value = 997
$ws.Range("K13").Value = 20000000
$ws.Range("L13").Value = 997
$ws.Range("M13").Value = -19999856
$ws.Range("N13").Value = -1285

$ws.Range("H98").Value = 355
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 355
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 355
$ws.Range("N98").Value = -6345

$ws.Range("H116").Value = 3779.6
$ws.Range("I116").Value = 2225
$ws.Range("J116").Value = 9998
$ws.Range("K116").Value = 2225
$ws.Range("L116").Value = 9998
$ws.Range("M116").Value = 69

$ws.Range("N133").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3779.6
$ws.Range("I3").Value = 2225
$ws.Range("J3").Value = 9998
$ws.Range("K3").Value = 2225
$ws.Range("L3").Value = 9998
$ws.Range("M3").Value = -2111

$ws.Range("H7").Value = 9500250
$ws.Range("I7").Value = 9500250
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 9500250
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -9500137

$ws.Range("H8").Value = 725
$ws.Range("I8").Value = 950
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 950
$ws.Range("L8").Value = 500
$ws.Range("M8").Value = -810
$ws.Range("N8").Value = -780

$ws.Range("N10").ClearContents()
$ws.Range("H10").Value = 602.5
$ws.Range("I10").Value = 602.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 602.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -462.5

$ws.Range("H11").Value = 2192.5
$ws.Range("I11").Value = 1980
$ws.Range("J11").Value = 2263.3333
$ws.Range("K11").Value = 1980
$ws.Range("L11").Value = 2263.3333
$ws.Range("M11").Value = -1840
$ws.Range("N11").Value = -2543.3333

$ws.Range("H12").Value = 893.5
$ws.Range("I12").Value = 893.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 893.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -725.5

$ws.Range("H20").Value = 2613.7778
$ws.Range("I20").Value = 2073.5715
$ws.Range("J20").Value = 4504.5
$ws.Range("K20").Value = 2073.5715
$ws.Range("L20").Value = 4504.5
$ws.Range("M20").Value = -1826.5715
$ws.Range("N20").Value = -4998.5

$ws.Range("H21").Value = 29527
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 29527
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 29527
$ws.Range("N21").Value = -29999

$ws.Range("H26").Value = 36616
$ws.Range("I26").Value = 36616
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 36616
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -36324

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 652.5
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 1005
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 1005
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -1231

$ws.Range("M3").ClearContents()
$ws.Range("H3").Value = 1003
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1003
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1003
$ws.Range("N3").Value = -1229

$ws.Range("H10").Value = 2206.2
$ws.Range("I10").Value = 257.75
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 257.75
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -118.75

$ws.Range("H11").Value = 5250
$ws.Range("I11").Value = 500
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 500
$ws.Range("L11").Value = 10000
$ws.Range("M11").Value = -360
$ws.Range("N11").Value = -10280

$ws.Range("H12").Value = 25326
$ws.Range("I12").Value = 652
$ws.Range("J12").Value = 50000
$ws.Range("K12").Value = 652
$ws.Range("L12").Value = 50000
$ws.Range("M12").Value = -482
$ws.Range("N12").Value = -50340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 933.3333
$ws.Range("I5").Value = 900.5
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 2701.5
$ws.Range("L5").Value = 2997
$ws.Range("M5").Value = -2589.5
$ws.Range("N5").Value = -3221

$ws.Range("H14").Value = 3332.3333
$ws.Range("I14").Value = 3332.3333
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 9996.999899999999
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -9823.999899999999

$ws.Range("H135").Value = 933.3333
$ws.Range("I135").Value = 900.5
$ws.Range("J135").Value = 999
$ws.Range("K135").Value = 8104.5
$ws.Range("L135").Value = 8991
$ws.Range("M135").Value = -5569.5
$ws.Range("N135").Value = -14061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1002300
$ws.Range("I3").Value = 1668666.6
$ws.Range("J3").Value = 2750
$ws.Range("K3").Value = 1668666.6
$ws.Range("L3").Value = 2750
$ws.Range("M3").Value = -1668550.6
$ws.Range("N3").Value = -2982

$ws.Range("H11").Value = 1773142.9
$ws.Range("I11").Value = 1235333.4
$ws.Range("J11").Value = 5000000
$ws.Range("K11").Value = 1235333.4
$ws.Range("L11").Value = 5000000
$ws.Range("M11").Value = -1235194.4
$ws.Range("N11").Value = -5000278

$ws.Range("N13").ClearContents()
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 136

$ws.Range("H98").Value = 16983.334
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 16983.334
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 16983.334
$ws.Range("N98").Value = -22973.334

$ws.Range("H107").Value = 3522.1177
$ws.Range("I107").Value = 501.41666
$ws.Range("J107").Value = 10771.8
$ws.Range("K107").Value = 501.41666
$ws.Range("L107").Value = 10771.8
$ws.Range("M107").Value = 1418.58334
$ws.Range("N107").Value = -14611.8

$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1994.3334
$ws.Range("I7").Value = 1994.3334
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1994.3334
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1882.3334

$ws.Range("H10").Value = 992.25
$ws.Range("I10").Value = 693
$ws.Range("J10").Value = 1890
$ws.Range("K10").Value = 693
$ws.Range("L10").Value = 1890
$ws.Range("M10").Value = -553
$ws.Range("N10").Value = -2170

$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0

$ws.Range("H99").Value = 19775.75
$ws.Range("I99").Value = 19775.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 19775.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -16780.75

$ws.Range("H100").Value = 1937.5
$ws.Range("I100").Value = 1937.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1937.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1396.5

$ws.Range("N112").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0

$ws.Range("H126").Value = 1994.3334
$ws.Range("I126").Value = 1994.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5983.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3513.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 85916.164
$ws.Range("I3").Value = 3099.4
$ws.Range("J3").Value = 500000
$ws.Range("K3").Value = 3099.4
$ws.Range("L3").Value = 500000
$ws.Range("M3").Value = -2985.4
$ws.Range("N3").Value = -500228

$ws.Range("H81").Value = 4366.1113
$ws.Range("I81").Value = 4286.875
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 8573.75
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -7512.75

$ws.Range("H84").Value = 4366.1113
$ws.Range("I84").Value = 4286.875
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 42868.75
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -37564.75

$ws.Range("H132").Value = 2164.8
$ws.Range("I132").Value = 1849.7142
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 5549.142599999999
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -3019.142599999999
